$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Insert a new heading paragraph "Datenbanken" in front of the
#    existing (empty) paragraph, and apply a "heading 1" style to it.
# ------------------------------------------------------------------
$firstPara = $d.Paragraphs(1)
$rng = $firstPara.Range
$rng.Collapse(1)                      # wdCollapseStart
$rng.InsertBefore("Datenbanken`r")

# ------------------------------------------------------------------
# 2. Create the paragraph style "berschrift1" (Überschrift 1 / Heading 1)
#    and its linked character style "berschrift1Zchn", matching the
#    definitions Word itself generates the first time "Heading 1" is
#    used in a German-language document.
# ------------------------------------------------------------------
$headingStyle = $d.Styles.Add("berschrift1", 1)          # wdStyleTypeParagraph
$headingStyle.NameLocal = "heading 1"
$headingStyle.BaseStyle = $d.Styles("Standard")
$headingStyle.NextParagraphStyle = $d.Styles("Standard")
$headingStyle.Priority = 9
$headingStyle.QuickStyle = $true

$headingCharStyle = $d.Styles.Add("berschrift1Zchn", 2)   # wdStyleTypeCharacter
$headingCharStyle.NameLocal = "Überschrift 1 Zchn"
$headingCharStyle.BaseStyle = $d.Styles("Absatz-Standardschriftart")
$headingCharStyle.Priority = 9

$headingStyle.LinkStyle = $headingCharStyle
$headingCharStyle.LinkStyle = $headingStyle

# Paragraph formatting for the heading style.
$pf = $headingStyle.ParagraphFormat
$pf.KeepWithNext = $true
$pf.KeepTogether = $true
$pf.SpaceBefore = 12
$pf.SpaceAfter = 0
$pf.OutlineLevel = 1

# Run formatting (font / color / size) shared by paragraph + character style.
# (matches the "majorHAnsi" theme font used by the real built-in Heading 1
#  style: ascii/hAnsi = the major-font's Latin typeface, i.e. "Calibri Light")
foreach ($st in @($headingStyle, $headingCharStyle)) {
    $fnt = $st.Font
    $fnt.Name = "Calibri Light"
    $fnt.Size = 16
    $fnt.SizeBi = 16
    $fnt.TextColor.ObjectThemeColor = 4   # msoThemeColorAccent1
}

# ------------------------------------------------------------------
# 3. Apply the heading style to the new first paragraph and set the text.
# ------------------------------------------------------------------
$firstPara2 = $d.Paragraphs(1)
$firstPara2.Style = $headingStyle
